$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tarantula")
$ws.Range("C2").Value = -1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 0
$ws.Range("I3").Value = 16
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0
$ws.Range("C5").Value = -1
$ws.Range("E5").Value = 0
$ws.Range("I5").Value = 6
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0
$ws.Range("I6").Value = 10
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 28
$ws.Range("C8").Value = -1
$ws.Range("E8").Value = 0
$ws.Range("I8").Value = 13
$ws.Range("C9").Value = -1
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 0
$ws.Range("I9").Value = 13
$ws.Range("C10").Value = -1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("I10").Value = 11
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 3
$ws.Range("C12").Value = -1
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 0
$ws.Range("I12").Value = 5
$ws.Range("C13").Value = -1
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 0
$ws.Range("C14").Value = -1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0
$ws.Range("C15").Value = -1
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("C16").Value = -1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("C17").Value = -1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 6
$ws.Range("C18").Value = -1
$ws.Range("E18").Value = 0
$ws.Range("C19").Value = -1
$ws.Range("E19").Value = 0
$ws.Range("I19").Value = 3
$ws.Range("C20").Value = -1
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 3
$ws.Range("C21").Value = -1
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 0
$ws.Range("I21").Value = 11
$ws.Range("C22").Value = -1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 2
$ws.Range("C23").Value = -1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("I23").Value = 4
$ws.Range("C24").Value = -1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 2
$ws.Range("I24").Value = 4
$ws.Range("C25").Value = -1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("I25").Value = 13
$ws.Range("C26").Value = -1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 3
$ws.Range("C27").Value = -1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 13
$ws.Range("C28").Value = -1
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 10
$ws.Range("C29").Value = -1
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 0
$ws.Range("I29").Value = 12
$ws.Range("C30").Value = -1
$ws.Range("E30").Value = 0
$ws.Range("I30").Value = 3
$ws.Range("C31").Value = -1
$ws.Range("E31").Value = 0
$ws.Range("I31").Value = 3
$ws.Range("C32").Value = -1
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 7
$ws.Range("C33").Value = -1
$ws.Range("E33").Value = 0
$ws.Range("C34").Value = -1
$ws.Range("D34").Value = 9
$ws.Range("E34").Value = 0
$ws.Range("I34").Value = 13

$ws = $wb.Worksheets.Item("Ochiai")
$ws.Range("C2").Value = -1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0
$ws.Range("I3").Value = 2
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0
$ws.Range("C5").Value = -1
$ws.Range("E5").Value = 0
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("C8").Value = -1
$ws.Range("E8").Value = 0
$ws.Range("C9").Value = -1
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 0
$ws.Range("I9").Value = 3
$ws.Range("C10").Value = -1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("I10").Value = 2
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 3
$ws.Range("C12").Value = -1
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 0
$ws.Range("I12").Value = 5
$ws.Range("C13").Value = -1
$ws.Range("E13").Value = 0
$ws.Range("C14").Value = -1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0
$ws.Range("C15").Value = -1
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("C16").Value = -1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("C17").Value = -1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 3
$ws.Range("C18").Value = -1
$ws.Range("E18").Value = 0
$ws.Range("C19").Value = -1
$ws.Range("E19").Value = 0
$ws.Range("C20").Value = -1
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 3
$ws.Range("C21").Value = -1
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 0
$ws.Range("I21").Value = 2
$ws.Range("C22").Value = -1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 14
$ws.Range("C23").Value = -1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("I23").Value = 3
$ws.Range("C24").Value = -1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 2
$ws.Range("I24").Value = 4
$ws.Range("C25").Value = -1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("I25").Value = 3
$ws.Range("C26").Value = -1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 2
$ws.Range("C27").Value = -1
$ws.Range("E27").Value = 0
$ws.Range("C28").Value = -1
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 9
$ws.Range("C29").Value = -1
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 0
$ws.Range("I29").Value = 2
$ws.Range("C30").Value = -1
$ws.Range("E30").Value = 0
$ws.Range("C31").Value = -1
$ws.Range("E31").Value = 0
$ws.Range("I31").Value = 2
$ws.Range("C32").Value = -1
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 4
$ws.Range("C33").Value = -1
$ws.Range("E33").Value = 0
$ws.Range("C34").Value = -1
$ws.Range("D34").Value = 9
$ws.Range("E34").Value = 0
$ws.Range("I34").Value = 3

$ws = $wb.Worksheets.Item("Op2")
$ws.Range("C2").Value = -1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0
$ws.Range("I3").Value = 2
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0
$ws.Range("C5").Value = -1
$ws.Range("E5").Value = 0
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("C8").Value = -1
$ws.Range("E8").Value = 0
$ws.Range("C9").Value = -1
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 0
$ws.Range("I9").Value = 3
$ws.Range("C10").Value = -1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("I10").Value = 2
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 3
$ws.Range("C12").Value = -1
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 0
$ws.Range("I12").Value = 5
$ws.Range("C13").Value = -1
$ws.Range("E13").Value = 0
$ws.Range("C14").Value = -1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0
$ws.Range("C15").Value = -1
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("C16").Value = -1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("C17").Value = -1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 3
$ws.Range("C18").Value = -1
$ws.Range("E18").Value = 0
$ws.Range("C19").Value = -1
$ws.Range("E19").Value = 0
$ws.Range("C20").Value = -1
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 3
$ws.Range("C21").Value = -1
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 0
$ws.Range("I21").Value = 2
$ws.Range("C22").Value = -1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 28
$ws.Range("C23").Value = -1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("I23").Value = 3
$ws.Range("C24").Value = -1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 2
$ws.Range("I24").Value = 4
$ws.Range("C25").Value = -1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("I25").Value = 3
$ws.Range("C26").Value = -1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 2
$ws.Range("C27").Value = -1
$ws.Range("E27").Value = 0
$ws.Range("C28").Value = -1
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 9
$ws.Range("C29").Value = -1
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 0
$ws.Range("I29").Value = 2
$ws.Range("C30").Value = -1
$ws.Range("E30").Value = 0
$ws.Range("C31").Value = -1
$ws.Range("E31").Value = 0
$ws.Range("I31").Value = 2
$ws.Range("C32").Value = -1
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 4
$ws.Range("C33").Value = -1
$ws.Range("E33").Value = 0
$ws.Range("C34").Value = -1
$ws.Range("D34").Value = 9
$ws.Range("E34").Value = 0
$ws.Range("I34").Value = 3

$ws = $wb.Worksheets.Item("Barinel")
$ws.Range("C2").Value = -1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 0
$ws.Range("I3").Value = 16
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0
$ws.Range("C5").Value = -1
$ws.Range("E5").Value = 0
$ws.Range("I5").Value = 6
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0
$ws.Range("I6").Value = 10
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 4
$ws.Range("I7").Value = 28
$ws.Range("C8").Value = -1
$ws.Range("E8").Value = 0
$ws.Range("I8").Value = 13
$ws.Range("C9").Value = -1
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 0
$ws.Range("I9").Value = 13
$ws.Range("C10").Value = -1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("I10").Value = 11
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 3
$ws.Range("C12").Value = -1
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 0
$ws.Range("I12").Value = 5
$ws.Range("C13").Value = -1
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 0
$ws.Range("C14").Value = -1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0
$ws.Range("C15").Value = -1
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("C16").Value = -1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("C17").Value = -1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 6
$ws.Range("C18").Value = -1
$ws.Range("E18").Value = 0
$ws.Range("C19").Value = -1
$ws.Range("E19").Value = 0
$ws.Range("I19").Value = 3
$ws.Range("C20").Value = -1
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 3
$ws.Range("C21").Value = -1
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 0
$ws.Range("I21").Value = 11
$ws.Range("C22").Value = -1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 2
$ws.Range("C23").Value = -1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("I23").Value = 4
$ws.Range("C24").Value = -1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 2
$ws.Range("I24").Value = 4
$ws.Range("C25").Value = -1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("I25").Value = 13
$ws.Range("C26").Value = -1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 3
$ws.Range("C27").Value = -1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 13
$ws.Range("C28").Value = -1
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 10
$ws.Range("C29").Value = -1
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 0
$ws.Range("I29").Value = 12
$ws.Range("C30").Value = -1
$ws.Range("E30").Value = 0
$ws.Range("I30").Value = 3
$ws.Range("C31").Value = -1
$ws.Range("E31").Value = 0
$ws.Range("I31").Value = 3
$ws.Range("C32").Value = -1
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 7
$ws.Range("C33").Value = -1
$ws.Range("E33").Value = 0
$ws.Range("C34").Value = -1
$ws.Range("D34").Value = 9
$ws.Range("E34").Value = 0
$ws.Range("I34").Value = 13

$ws = $wb.Worksheets.Item("Dstar")
$ws.Range("C2").Value = -1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0
$ws.Range("I3").Value = 2
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0
$ws.Range("C5").Value = -1
$ws.Range("E5").Value = 0
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("C8").Value = -1
$ws.Range("E8").Value = 0
$ws.Range("C9").Value = -1
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 0
$ws.Range("I9").Value = 3
$ws.Range("C10").Value = -1
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("I10").Value = 2
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 3
$ws.Range("C12").Value = -1
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 0
$ws.Range("I12").Value = 5
$ws.Range("C13").Value = -1
$ws.Range("E13").Value = 0
$ws.Range("C14").Value = -1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0
$ws.Range("C15").Value = -1
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("C16").Value = -1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("C17").Value = -1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 3
$ws.Range("C18").Value = -1
$ws.Range("E18").Value = 0
$ws.Range("C19").Value = -1
$ws.Range("E19").Value = 0
$ws.Range("C20").Value = -1
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 3
$ws.Range("C21").Value = -1
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 0
$ws.Range("I21").Value = 2
$ws.Range("C22").Value = -1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 14
$ws.Range("C23").Value = -1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("I23").Value = 3
$ws.Range("C24").Value = -1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 2
$ws.Range("I24").Value = 4
$ws.Range("C25").Value = -1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("I25").Value = 3
$ws.Range("C26").Value = -1
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
$ws.Range("C27").Value = -1
$ws.Range("E27").Value = 0
$ws.Range("C28").Value = -1
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 9
$ws.Range("C29").Value = -1
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 0
$ws.Range("I29").Value = 2
$ws.Range("C30").Value = -1
$ws.Range("E30").Value = 0
$ws.Range("C31").Value = -1
$ws.Range("E31").Value = 0
$ws.Range("I31").Value = 2
$ws.Range("C32").Value = -1
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 4
$ws.Range("C33").Value = -1
$ws.Range("E33").Value = 0
$ws.Range("C34").Value = -1
$ws.Range("D34").Value = 9
$ws.Range("E34").Value = 0
$ws.Range("I34").Value = 3
